$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-04 Sunday" "2024-08-05 Monday"

Replace-Text "865×6=5190" "939×7=6573"
Replace-Text "253×9=2277" "679×3=2037"
Replace-Text "827×9=7443" "484×9=4356"
Replace-Text "391×5=1955" "539×8=4312"
Replace-Text "742×7=5194" "437×7=3059"

Replace-Text "533×3=1599" "690×2=1380"
Replace-Text "804×8=6432" "593×6=3558"
Replace-Text "946×7=6622" "556×8=4448"
Replace-Text "357×2=714" "651×7=4557"
Replace-Text "931×5=4655" "215×6=1290"

Replace-Text "328×8=2624" "170×3=510"
Replace-Text "398×4=1592" "199×6=1194"
Replace-Text "845×8=6760" "294×3=882"
Replace-Text "730×2=1460" "506×9=4554"
Replace-Text "251×8=2008" "490×4=1960"

Replace-Text "890×8=7120" "231×3=693"
Replace-Text "357×3=1071" "359×2=718"
Replace-Text "204×9=1836" "770×3=2310"
Replace-Text "746×2=1492" "531×9=4779"
Replace-Text "590×5=2950" "347×4=1388"

Replace-Text "813×3=2439" "264×3=792"
Replace-Text "634×6=3804" "311×5=1555"
Replace-Text "508×7=3556" "558×7=3906"
Replace-Text "626×6=3756" "589×7=4123"
Replace-Text "469×4=1876" "951×3=2853"
